$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "ODI Batting Extra" after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- Header row values ---
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Reuse the existing bold/bordered/centered header style (same as the other
# sheets' header rows) by copying the format from "ODI Bowling"!A1 instead of
# re-declaring Font/Border/Alignment (which would mint brand-new style
# records).
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- A temporary scratch sheet is used to stage text-typed values so that
# numeric-looking strings ("3955", "0", "0.37%", ...) land in the target
# cells as plain text (matching the source data) instead of being
# auto-coerced to numbers, all without leaving any extra number-format /
# style behind on the target cells. Copy + PasteSpecial(xlPasteValues)
# carries just the text-ness of the source cell, not its style. ---
$scratchWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$scratchWs.Name = "ScratchTemp"
$scratch = $scratchWs.Range("A1")
$scratch.NumberFormat = "@"

function Set-TextValue($targetCell, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)
}

# --- Row 2 ---
Set-TextValue $ws.Range("A2") "3955"
$ws.Range("B2").Value = 8
Set-TextValue $ws.Range("C2") "0"
Set-TextValue $ws.Range("D2") "0"
Set-TextValue $ws.Range("E2") "0.37%"
Set-TextValue $ws.Range("F2") "NO"

# --- Row 3 ---
Set-TextValue $ws.Range("A3") "4529"
$ws.Range("B3").Value = 8
Set-TextValue $ws.Range("C3") "0"
Set-TextValue $ws.Range("D3") "0"
Set-TextValue $ws.Range("E3") "0.71%"
Set-TextValue $ws.Range("F3") "NO"

# Remove the scratch sheet — it was only a staging area.
[void]$scratchWs.Delete()

[void]$ws.Range("A1").Select()
